$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value2 = 4491.6665
$ws.Range("I132").Value2 = 0
$ws.Range("J132").Value2 = 4491.6665
$ws.Range("K132").Value2 = 0
$ws.Range("L132").Value2 = 13474.9995
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value2 = -18534.9995

$ws.Range("H134").Value2 = 30780
$ws.Range("J134").Value2 = 30780
$ws.Range("L134").Value2 = 30780
$ws.Range("N134").Value2 = -40920

$ws.Range("H135").Value2 = 885.1064
$ws.Range("I135").Value2 = 648.7805
$ws.Range("J135").Value2 = 2500
$ws.Range("K135").Value2 = 5839.0245
$ws.Range("L135").Value2 = 22500
$ws.Range("M135").Value2 = -3304.0245
$ws.Range("N135").Value2 = -27570

$ws.Range("H137").Value2 = 3561.652
$ws.Range("I137").Value2 = 3999.9375
$ws.Range("J137").Value2 = 2559.8572
$ws.Range("K137").Value2 = 11999.8125
$ws.Range("L137").Value2 = 7679.571599999999
$ws.Range("M137").Value2 = -9449.8125
$ws.Range("N137").Value2 = -12779.5716

$ws.Range("H138").Value2 = 5362.92
$ws.Range("I138").Value2 = 1427.5
$ws.Range("J138").Value2 = 9626.291999999999
$ws.Range("K138").Value2 = 4282.5
$ws.Range("L138").Value2 = 28878.876
$ws.Range("M138").Value2 = 857.5
$ws.Range("N138").Value2 = -39158.876

$ws.Range("H141").Value2 = 697917.4399999999
$ws.Range("I141").Value2 = 6040.84
$ws.Range("J141").Value2 = 2139327
$ws.Range("K141").Value2 = 18122.52
$ws.Range("L141").Value2 = 6417981
$ws.Range("M141").Value2 = -12942.52
$ws.Range("N141").Value2 = -6428341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 8930278
$ws.Range("I2").Value2 = 19231990
$ws.Range("K2").Value2 = 19231990
$ws.Range("M2").Value2 = -19231877

$ws.Range("H61").Value2 = 1298.2683
$ws.Range("I61").Value2 = 681.86487
$ws.Range("J61").Value2 = 7000
$ws.Range("K61").Value2 = 681.86487
$ws.Range("L61").Value2 = 7000
$ws.Range("M61").Value2 = -469.86487
$ws.Range("N61").Value2 = -7424

$ws.Range("H74").Value2 = 1204.421
$ws.Range("I74").Value2 = 1024
$ws.Range("J74").Value2 = 2166.6667
$ws.Range("K74").Value2 = 1024
$ws.Range("L74").Value2 = 2166.6667
$ws.Range("M74").Value2 = -150
$ws.Range("N74").Value2 = -3914.6667

$ws.Range("H77").Value2 = 1204.421
$ws.Range("I77").Value2 = 1024
$ws.Range("J77").Value2 = 2166.6667
$ws.Range("K77").Value2 = 5120
$ws.Range("L77").Value2 = 10833.3335
$ws.Range("M77").Value2 = -752
$ws.Range("N77").Value2 = -19569.3335

$ws.Range("H97").Value2 = 689.5263
$ws.Range("I97").Value2 = 396.06668
$ws.Range("J97").Value2 = 1790
$ws.Range("K97").Value2 = 396.06668
$ws.Range("L97").Value2 = 1790
$ws.Range("M97").Value2 = 99.93331999999998
$ws.Range("N97").Value2 = -2782

$ws.Range("H102").Value2 = 2111.4285
$ws.Range("I102").Value2 = 1523.6364
$ws.Range("J102").Value2 = 4266.6665
$ws.Range("K102").Value2 = 1523.6364
$ws.Range("L102").Value2 = 4266.6665
$ws.Range("M102").Value2 = 98.36359999999991
$ws.Range("N102").Value2 = -7510.6665

$ws.Range("H116").Value2 = 8930278
$ws.Range("I116").Value2 = 19231990
$ws.Range("K116").Value2 = 19231990
$ws.Range("M116").Value2 = -19229696

$ws.Range("H122").Value2 = 1668.1041
$ws.Range("I122").Value2 = 1239.8649
$ws.Range("J122").Value2 = 3108.5454
$ws.Range("K122").Value2 = 3719.5947
$ws.Range("L122").Value2 = 9325.636200000001
$ws.Range("M122").Value2 = -1269.5947
$ws.Range("N122").Value2 = -14225.6362

$ws.Range("H134").Value2 = 84900
$ws.Range("J134").Value2 = 84900
$ws.Range("L134").Value2 = 84900
$ws.Range("N134").Value2 = -95040

$ws.Range("H136").Value2 = 1298.2683
$ws.Range("I136").Value2 = 681.86487
$ws.Range("J136").Value2 = 7000
$ws.Range("K136").Value2 = 2045.59461
$ws.Range("L136").Value2 = 21000
$ws.Range("M136").Value2 = 504.4053899999999
$ws.Range("N136").Value2 = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 8930278
$ws.Range("I3").Value2 = 19231990
$ws.Range("K3").Value2 = 19231990
$ws.Range("M3").Value2 = -19231876

$ws.Range("H86").Value2 = 613977.8
$ws.Range("I86").Value2 = 863342.5
$ws.Range("J86").Value2 = 59834.11
$ws.Range("K86").Value2 = 863342.5
$ws.Range("L86").Value2 = 59834.11
$ws.Range("M86").Value2 = -862219.5
$ws.Range("N86").Value2 = -62080.11

$ws.Range("H89").Value2 = 613977.8
$ws.Range("I89").Value2 = 863342.5
$ws.Range("J89").Value2 = 59834.11
$ws.Range("K89").Value2 = 4316712.5
$ws.Range("L89").Value2 = 299170.55
$ws.Range("M89").Value2 = -4311096.5
$ws.Range("N89").Value2 = -310402.55

$ws.Range("H107").Value2 = 2006
$ws.Range("I107").Value2 = 1605.5
$ws.Range("J107").Value2 = 3183.9412
$ws.Range("K107").Value2 = 1605.5
$ws.Range("L107").Value2 = 3183.9412
$ws.Range("M107").Value2 = 314.5
$ws.Range("N107").Value2 = -7023.9412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 3029.0444
$ws.Range("I31").Value2 = 1866.909
$ws.Range("J31").Value2 = 6224.9165
$ws.Range("K31").Value2 = 1866.909
$ws.Range("L31").Value2 = 6224.9165
$ws.Range("M31").Value2 = -1571.909
$ws.Range("N31").Value2 = -6814.9165

$ws.Range("H34").Value2 = 3029.0444
$ws.Range("I34").Value2 = 1866.909
$ws.Range("J34").Value2 = 6224.9165
$ws.Range("K34").Value2 = 1866.909
$ws.Range("L34").Value2 = 6224.9165
$ws.Range("M34").Value2 = -1664.909
$ws.Range("N34").Value2 = -6628.9165

$ws.Range("H58").Value2 = 7044655.5
$ws.Range("I58").Value2 = 1461.6666
$ws.Range("J58").Value2 = 29417152
$ws.Range("K58").Value2 = 1461.6666
$ws.Range("L58").Value2 = 29417152
$ws.Range("M58").Value2 = -1258.6666
$ws.Range("N58").Value2 = -29417558

$ws.Range("H99").Value2 = 1784.4
$ws.Range("I99").Value2 = 1244.5333
$ws.Range("J99").Value2 = 3404
$ws.Range("K99").Value2 = 1244.5333
$ws.Range("L99").Value2 = 3404
$ws.Range("M99").Value2 = 253.4666999999999
$ws.Range("N99").Value2 = -6400

$ws.Range("H107").Value2 = 1877.4736
$ws.Range("I107").Value2 = 535.9
$ws.Range("J107").Value2 = 3368.111
$ws.Range("K107").Value2 = 535.9
$ws.Range("L107").Value2 = 3368.111
$ws.Range("M107").Value2 = 1384.1
$ws.Range("N107").Value2 = -7208.111

$ws.Range("H111").Value2 = 0
$ws.Range("J111").Value2 = 0
$ws.Range("L111").Value2 = 0
$ws.Range("N111").ClearContents()

$ws.Range("H126").Value2 = 1784.4
$ws.Range("I126").Value2 = 1244.5333
$ws.Range("J126").Value2 = 3404
$ws.Range("K126").Value2 = 3733.5999
$ws.Range("L126").Value2 = 10212
$ws.Range("M126").Value2 = -1263.5999
$ws.Range("N126").Value2 = -15152

$ws.Range("H136").Value2 = 7044655.5
$ws.Range("I136").Value2 = 1461.6666
$ws.Range("J136").Value2 = 29417152
$ws.Range("K136").Value2 = 4384.9998
$ws.Range("L136").Value2 = 88251456
$ws.Range("M136").Value2 = -1834.9998
$ws.Range("N136").Value2 = -88256556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value2 = 2212.5
$ws.Range("I59").Value2 = 50
$ws.Range("J59").Value2 = 2933.3333
$ws.Range("K59").Value2 = 150
$ws.Range("L59").Value2 = 8799.999899999999
$ws.Range("M59").Value2 = 390
$ws.Range("N59").Value2 = -9879.999899999999

$ws.Range("H95").Value2 = 4000
$ws.Range("J95").Value2 = 4000
$ws.Range("L95").Value2 = 12000
$ws.Range("N95").Value2 = -16118

$ws.Range("H131").Value2 = 2527
$ws.Range("I131").Value2 = 2584
$ws.Range("J131").Value2 = 2470
$ws.Range("K131").Value2 = 7752
$ws.Range("L131").Value2 = 7410
$ws.Range("M131").Value2 = -2712
$ws.Range("N131").Value2 = -17490

$ws.Range("H133").Value2 = 4526.4375
$ws.Range("I133").Value2 = 5827.143
$ws.Range("J133").Value2 = 3514.7778
$ws.Range("K133").Value2 = 17481.429
$ws.Range("L133").Value2 = 10544.3334
$ws.Range("M133").Value2 = -12421.429
$ws.Range("N133").Value2 = -20664.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2436.875
$ws.Range("I80").Value2 = 2249.1667
$ws.Range("J80").Value2 = 3000
$ws.Range("K80").Value2 = 2249.1667
$ws.Range("L80").Value2 = 3000
$ws.Range("M80").Value2 = -1251.1667
$ws.Range("N80").Value2 = -4996

$ws.Range("H83").Value2 = 2436.875
$ws.Range("I83").Value2 = 2249.1667
$ws.Range("J83").Value2 = 3000
$ws.Range("K83").Value2 = 11245.8335
$ws.Range("L83").Value2 = 15000
$ws.Range("M83").Value2 = -6253.833500000001
$ws.Range("N83").Value2 = -24984

$ws.Range("H126").Value2 = 4390.909
$ws.Range("I126").Value2 = 3760
$ws.Range("J126").Value2 = 4916.6665
$ws.Range("K126").Value2 = 11280
$ws.Range("L126").Value2 = 14749.9995
$ws.Range("M126").Value2 = -8810
$ws.Range("N126").Value2 = -19689.9995

$ws.Range("H132").Value2 = 2552.7017
$ws.Range("I132").Value2 = 2166.4888
$ws.Range("J132").Value2 = 4001
$ws.Range("K132").Value2 = 6499.4664
$ws.Range("L132").Value2 = 12003
$ws.Range("M132").Value2 = -3969.4664
$ws.Range("N132").Value2 = -17063

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 125001580
$ws.Range("I22").Value2 = 200000340
$ws.Range("J22").Value2 = 3660
$ws.Range("K22").Value2 = 200000340
$ws.Range("L22").Value2 = 3660
$ws.Range("M22").Value2 = -200000045
$ws.Range("N22").Value2 = -4250

$ws.Range("H27").Value2 = 125001580
$ws.Range("I27").Value2 = 200000340
$ws.Range("J27").Value2 = 3660
$ws.Range("K27").Value2 = 200000340
$ws.Range("L27").Value2 = 3660
$ws.Range("M27").Value2 = -200000233
$ws.Range("N27").Value2 = -3874

$ws.Range("H132").Value2 = 2051.037
$ws.Range("I132").Value2 = 1309.2565
$ws.Range("J132").Value2 = 3979.6667
$ws.Range("K132").Value2 = 3927.7695
$ws.Range("L132").Value2 = 11939.0001
$ws.Range("M132").Value2 = -1397.7695
$ws.Range("N132").Value2 = -16999.0001

$ws.Range("H136").Value2 = 1524.0364
$ws.Range("I136").Value2 = 1048.2916
$ws.Range("J136").Value2 = 4786.2856
$ws.Range("K136").Value2 = 3144.8748
$ws.Range("L136").Value2 = 14358.8568
$ws.Range("M136").Value2 = -594.8748000000001
$ws.Range("N136").Value2 = -19458.8568
